$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.699.34"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "3.566.48"

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'607.54"
$ws.Range("E5").Value = "  -0.12%  "

$ws.Range("D6").Value = "'145.47"
$ws.Range("E6").Value = "  +1.06%  "

$ws.Range("D7").Value = "3.565.57"
$ws.Range("E7").Value = "  +0.83%  "

$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("E9").Value = "  +3.63%  "

$ws.Range("E10").Value = "  +0.02%  "

$ws.Range("D11").Value = "'7.96"
$ws.Range("E11").Value = "  -1.60%  "

$ws.Range("E12").Value = "  +0.88%  "

$ws.Range("D13").Value = "4.170.22"
$ws.Range("E13").Value = "  +0.80%  "

$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("D15").Value = "'30.13"

$ws.Range("D16").Value = "3.545.12"
$ws.Range("E16").Value = "  +0.36%  "

$ws.Range("D17").Value = "66.714.99"
$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("E18").Value = "  +0.32%  "

$ws.Range("D19").Value = "'11.54"
$ws.Range("E19").Value = "  +5.44%  "

$ws.Range("D20").Value = "'6.20"
$ws.Range("E20").Value = "  -0.17%  "

$ws.Range("D21").Value = "'14.90"
$ws.Range("E21").Value = "  -0.42%  "

$ws.Range("D22").Value = "'432.67"
$ws.Range("E22").Value = "  +1.58%  "

$ws.Range("D23").Value = "'0.612"
$ws.Range("E23").Value = "  +1.72%  "

$ws.Range("D24").Value = "'79.93"
$ws.Range("E24").Value = "  +1.57%  "

$ws.Range("D25").Value = "3.707.26"
$ws.Range("E25").Value = "  +0.76%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("D27").Value = "'0.0000121"
$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").Value = "'8.07"
$ws.Range("E28").Value = "  -0.41%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'9.22"
$ws.Range("E29").Value = "  +0.26%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.51"
$ws.Range("E30").Value = "  +1.36%  "

$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("E32").Value = "  -1.71%  "

$ws.Range("D33").Value = "3.560.67"

$ws.Range("D34").Value = "'25.40"

$ws.Range("E35").Value = "  -3.61%  "

$ws.Range("D36").Value = "'7.89"
$ws.Range("E36").Value = "  +0.83%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").Value = "'1.73"
$ws.Range("E38").Value = "  -1.83%  "

$ws.Range("D39").Value = "'5.63"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("D40").Value = "'174.67"
$ws.Range("E40").Value = "  +1.26%  "

$ws.Range("D41").Value = "'0.0852"
$ws.Range("E41").Value = "  -0.57%  "

$ws.Range("D42").Value = "'5.21"
$ws.Range("E42").Value = "  +0.51%  "

$ws.Range("D43").Value = "'0.889"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("D44").Value = "'1.96"
$ws.Range("E44").Value = "  +3.19%  "

$ws.Range("D45").Value = "'46.15"
$ws.Range("E45").Value = "  +1.37%  "

$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").Value = "'2.54"
$ws.Range("E47").Value = "  +5.66%  "

$ws.Range("E48").Value = "  -2.35%  "

$ws.Range("D49").Value = "'25.19"
$ws.Range("E49").Value = "  -3.46%  "

$ws.Range("D50").Value = "'23.63"
$ws.Range("E50").Value = "  +4.38%  "

$ws.Range("D51").Value = "'7.17"
$ws.Range("E51").Value = "  +0.41%  "
